$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15016.5
$ws.Range("I18").Value = 16419.8
$ws.Range("K18").Value = 16419.8
$ws.Range("M18").Value = -16135.8

$ws.Range("H39").Value = 41.736843
$ws.Range("I39").Value = 41.285713
$ws.Range("J39").Value = 43
$ws.Range("K39").Value = 123.857139
$ws.Range("L39").Value = 129
$ws.Range("M39").Value = 172.142861
$ws.Range("N39").Value = -721

$ws.Range("H40").Value = 4124.6313
$ws.Range("I40").Value = 2966.6667
$ws.Range("K40").Value = 2966.6667
$ws.Range("M40").Value = -2791.6667

$ws.Range("H80").Value = 698.069
$ws.Range("I80").Value = 427.75
$ws.Range("J80").Value = 888.8823
$ws.Range("K80").Value = 1283.25
$ws.Range("L80").Value = 2666.6469
$ws.Range("M80").Value = -285.25
$ws.Range("N80").Value = -4662.6469

$ws.Range("H83").Value = 698.069
$ws.Range("I83").Value = 427.75
$ws.Range("J83").Value = 888.8823
$ws.Range("K83").Value = 3849.75
$ws.Range("L83").Value = 7999.9407
$ws.Range("M83").Value = 1142.25
$ws.Range("N83").Value = -17983.9407

$ws.Range("H93").Value = 90000
$ws.Range("J93").Value = 90000
$ws.Range("L93").Value = 90000
$ws.Range("N93").Value = -94992

$ws.Range("H107").Value = 460.95
$ws.Range("I107").Value = 502.83334
$ws.Range("J107").Value = 84
$ws.Range("K107").Value = 502.83334
$ws.Range("L107").Value = 84
$ws.Range("M107").Value = 1417.16666
$ws.Range("N107").Value = -3924

$ws.Range("H118").Value = 1468.0625
$ws.Range("I118").Value = 1198.7778
$ws.Range("J118").Value = 1814.2858
$ws.Range("K118").Value = 3596.3334
$ws.Range("L118").Value = 5442.857400000001
$ws.Range("M118").Value = -1939.3334
$ws.Range("N118").Value = -8756.857400000001

$ws.Range("H137").Value = 1521.4584
$ws.Range("I137").Value = 1531.9546
$ws.Range("J137").Value = 1406
$ws.Range("K137").Value = 4595.8638
$ws.Range("L137").Value = 4218
$ws.Range("M137").Value = -2045.8638
$ws.Range("N137").Value = -9318

$ws.Range("H139").Value = 70416.664
$ws.Range("J139").Value = 70416.664
$ws.Range("L139").Value = 70416.664
$ws.Range("N139").Value = -80696.664

$ws.Range("H141").Value = 4510.222
$ws.Range("I141").Value = 4047.5
$ws.Range("J141").Value = 4642.4287
$ws.Range("K141").Value = 12142.5
$ws.Range("L141").Value = 13927.2861
$ws.Range("M141").Value = -6962.5
$ws.Range("N141").Value = -24287.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4067
$ws.Range("I61").Value = 4067
$ws.Range("K61").Value = 4067
$ws.Range("M61").Value = -3855

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H136").Value = 4067
$ws.Range("I136").Value = 4067
$ws.Range("K136").Value = 12201
$ws.Range("M136").Value = -9651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1033.8636
$ws.Range("I20").Value = 1060.5
$ws.Range("K20").Value = 1060.5
$ws.Range("M20").Value = -813.5

$ws.Range("H99").Value = 2270.842
$ws.Range("I99").Value = 1171
$ws.Range("J99").Value = 2912.4167
$ws.Range("K99").Value = 1171
$ws.Range("L99").Value = 2912.4167
$ws.Range("M99").Value = 327
$ws.Range("N99").Value = -5908.4167

$ws.Range("H107").Value = 4012.9473
$ws.Range("I107").Value = 1880.75
$ws.Range("K107").Value = 1880.75
$ws.Range("M107").Value = 39.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1472.6
$ws.Range("I22").Value = 1358.1111
$ws.Range("J22").Value = 1644.3334
$ws.Range("K22").Value = 1358.1111
$ws.Range("L22").Value = 1644.3334
$ws.Range("M22").Value = -1008.1111
$ws.Range("N22").Value = -2344.3334

$ws.Range("H31").Value = 3123.7368
$ws.Range("I31").Value = 1532.1818
$ws.Range("K31").Value = 1532.1818
$ws.Range("M31").Value = -1237.1818

$ws.Range("H34").Value = 3123.7368
$ws.Range("I34").Value = 1532.1818
$ws.Range("K34").Value = 1532.1818
$ws.Range("M34").Value = -1330.1818

$ws.Range("H45").Value = 20000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H132").Value = 1820.3334
$ws.Range("I132").Value = 1584.7333
$ws.Range("K132").Value = 4754.199900000001
$ws.Range("M132").Value = -2224.199900000001

$ws.Range("H134").Value = 5130.857
$ws.Range("I134").Value = 3666.625
$ws.Range("K134").Value = 10999.875
$ws.Range("M134").Value = -8464.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1074.75
$ws.Range("I64").Value = 1074.75
$ws.Range("K64").Value = 3224.25
$ws.Range("M64").Value = -2954.25

$ws.Range("H67").Value = 1074.75
$ws.Range("I67").Value = 1074.75
$ws.Range("K67").Value = 3224.25
$ws.Range("M67").Value = -2288.25

$ws.Range("H97").Value = 921.8125
$ws.Range("I97").Value = 271.5
$ws.Range("J97").Value = 1138.5834
$ws.Range("K97").Value = 814.5
$ws.Range("L97").Value = 3415.7502
$ws.Range("M97").Value = -318.5
$ws.Range("N97").Value = -4407.7502

$ws.Range("H139").Value = 6318.6665
$ws.Range("J139").Value = 4895.8335
$ws.Range("L139").Value = 14687.5005
$ws.Range("N139").Value = -24967.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1542.8182
$ws.Range("I107").Value = 482
$ws.Range("J107").Value = 1940.625
$ws.Range("K107").Value = 482
$ws.Range("L107").Value = 1940.625
$ws.Range("M107").Value = 1438
$ws.Range("N107").Value = -5780.625

$ws.Range("H132").Value = 2021.5
$ws.Range("I132").Value = 1925.683
$ws.Range("J132").Value = 3331
$ws.Range("K132").Value = 5777.049
$ws.Range("L132").Value = 9993
$ws.Range("M132").Value = -3247.049
$ws.Range("N132").Value = -15053

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4197.421
$ws.Range("I7").Value = 3439.0715
$ws.Range("K7").Value = 3439.0715
$ws.Range("M7").Value = -3327.0715

$ws.Range("H126").Value = 4197.421
$ws.Range("I126").Value = 3439.0715
$ws.Range("K126").Value = 10317.2145
$ws.Range("M126").Value = -7847.2145

$ws.Range("H132").Value = 3808.0588
$ws.Range("I132").Value = 3608.8
$ws.Range("K132").Value = 10826.4
$ws.Range("M132").Value = -8296.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 122249.5
$ws.Range("J16").Value = 122249.5
$ws.Range("L16").Value = 122249.5
